$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 48418.184
$ws.Range("I64").Value = 65437.5
$ws.Range("J64").Value = 3033.3333
$ws.Range("K64").Value = 65437.5
$ws.Range("L64").Value = 3033.3333
$ws.Range("M64").Value = -65189.5
$ws.Range("N64").Value = -3529.3333
$ws.Range("H67").Value = 48418.184
$ws.Range("I67").Value = 65437.5
$ws.Range("J67").Value = 3033.3333
$ws.Range("K67").Value = 65437.5
$ws.Range("L67").Value = 3033.3333
$ws.Range("M67").Value = -64579.5
$ws.Range("N67").Value = -4749.3333
$ws.Range("H74").Value = 3066.2222
$ws.Range("I74").Value = 3062
$ws.Range("J74").Value = 3100
$ws.Range("K74").Value = 3062
$ws.Range("L74").Value = 3100
$ws.Range("M74").Value = -2126
$ws.Range("N74").Value = -4972
$ws.Range("H77").Value = 3066.2222
$ws.Range("I77").Value = 3062
$ws.Range("J77").Value = 3100
$ws.Range("K77").Value = 15310
$ws.Range("L77").Value = 15500
$ws.Range("M77").Value = -10630
$ws.Range("N77").Value = -24860
$ws.Range("H100").Value = 989.4545000000001
$ws.Range("I100").Value = 771.75
$ws.Range("J100").Value = 1113.8572
$ws.Range("K100").Value = 771.75
$ws.Range("L100").Value = 1113.8572
$ws.Range("M100").Value = -230.75
$ws.Range("N100").Value = -2195.8572
$ws.Range("H129").Value = 3622.842
$ws.Range("J129").Value = 1119.0344
$ws.Range("L129").Value = 3357.1032
$ws.Range("N129").Value = -13357.1032

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29005.838
$ws.Range("I32").Value = 4675.385
$ws.Range("J32").Value = 155524.2
$ws.Range("K32").Value = 4675.385
$ws.Range("L32").Value = 155524.2
$ws.Range("M32").Value = -4388.385
$ws.Range("N32").Value = -156098.2
$ws.Range("H97").Value = 29868.885
$ws.Range("I97").Value = 34217
$ws.Range("J97").Value = 3780.2
$ws.Range("K97").Value = 34217
$ws.Range("L97").Value = 3780.2
$ws.Range("M97").Value = -33721
$ws.Range("N97").Value = -4772.2
$ws.Range("H102").Value = 65138.562
$ws.Range("I102").Value = 201792
$ws.Range("J102").Value = 3023.3635
$ws.Range("K102").Value = 201792
$ws.Range("L102").Value = 3023.3635
$ws.Range("M102").Value = -200170
$ws.Range("N102").Value = -6267.363499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 54744.105
$ws.Range("I20").Value = 64838.875
$ws.Range("J20").Value = 905.3333
$ws.Range("K20").Value = 64838.875
$ws.Range("L20").Value = 905.3333
$ws.Range("M20").Value = -64591.875
$ws.Range("N20").Value = -1399.3333
$ws.Range("H86").Value = 41996.645
$ws.Range("I86").Value = 86669.62
$ws.Range("J86").Value = 3280.0667
$ws.Range("K86").Value = 86669.62
$ws.Range("L86").Value = 3280.0667
$ws.Range("M86").Value = -85546.62
$ws.Range("N86").Value = -5526.066699999999
$ws.Range("H89").Value = 41996.645
$ws.Range("I89").Value = 86669.62
$ws.Range("J89").Value = 3280.0667
$ws.Range("K89").Value = 433348.1
$ws.Range("L89").Value = 16400.3335
$ws.Range("M89").Value = -427732.1
$ws.Range("N89").Value = -27632.3335
$ws.Range("H99").Value = 1483.1562
$ws.Range("I99").Value = 1307.4546
$ws.Range("K99").Value = 1307.4546
$ws.Range("M99").Value = 190.5454
$ws.Range("H105").Value = 57688.89
$ws.Range("I105").Value = 45735.22
$ws.Range("J105").Value = 78837.69500000001
$ws.Range("K105").Value = 45735.22
$ws.Range("L105").Value = 78837.69500000001
$ws.Range("M105").Value = -43988.22
$ws.Range("N105").Value = -82331.69500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24468.885
$ws.Range("I31").Value = 860.40625
$ws.Range("J31").Value = 44887.027
$ws.Range("K31").Value = 860.40625
$ws.Range("L31").Value = 44887.027
$ws.Range("M31").Value = -565.40625
$ws.Range("N31").Value = -45477.027
$ws.Range("H34").Value = 24468.885
$ws.Range("I34").Value = 860.40625
$ws.Range("J34").Value = 44887.027
$ws.Range("K34").Value = 860.40625
$ws.Range("L34").Value = 44887.027
$ws.Range("M34").Value = -658.40625
$ws.Range("N34").Value = -45291.027
$ws.Range("H41").Value = 12805
$ws.Range("J41").Value = 14356
$ws.Range("L41").Value = 14356
$ws.Range("N41").Value = -15212
$ws.Range("H62").Value = 2662.5
$ws.Range("J62").Value = 2662.5
$ws.Range("L62").Value = 2662.5
$ws.Range("N62").Value = -3910.5
$ws.Range("H65").Value = 2662.5
$ws.Range("J65").Value = 2662.5
$ws.Range("L65").Value = 13312.5
$ws.Range("N65").Value = -19552.5
$ws.Range("H105").Value = 2727.9
$ws.Range("I105").Value = 3009.7144
$ws.Range("J105").Value = 2070.3333
$ws.Range("K105").Value = 3009.7144
$ws.Range("L105").Value = 2070.3333
$ws.Range("M105").Value = -1262.7144
$ws.Range("N105").Value = -5564.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 78337.81
$ws.Range("I70").Value = 146632.58
$ws.Range("J70").Value = 4789.615
$ws.Range("K70").Value = 146632.58
$ws.Range("L70").Value = 4789.615
$ws.Range("M70").Value = -146362.58
$ws.Range("N70").Value = -5329.615
$ws.Range("H73").Value = 78337.81
$ws.Range("I73").Value = 146632.58
$ws.Range("J73").Value = 4789.615
$ws.Range("K73").Value = 146632.58
$ws.Range("L73").Value = 4789.615
$ws.Range("M73").Value = -145696.58
$ws.Range("N73").Value = -6661.615
$ws.Range("H80").Value = 111227190
$ws.Range("I80").Value = 250258750
$ws.Range("J80").Value = 1950
$ws.Range("K80").Value = 250258750
$ws.Range("L80").Value = 1950
$ws.Range("M80").Value = -250257752
$ws.Range("N80").Value = -3946
$ws.Range("H83").Value = 111227190
$ws.Range("I83").Value = 250258750
$ws.Range("J83").Value = 1950
$ws.Range("K83").Value = 1251293750
$ws.Range("L83").Value = 9750
$ws.Range("M83").Value = -1251288758
$ws.Range("N83").Value = -19734
$ws.Range("H97").Value = 76924500
$ws.Range("I97").Value = 100001700
$ws.Range("J97").Value = 470.66666
$ws.Range("K97").Value = 100001700
$ws.Range("L97").Value = 470.66666
$ws.Range("M97").Value = -100001204
$ws.Range("N97").Value = -1462.66666
$ws.Range("H126").Value = 6539870
$ws.Range("I126").Value = 4667.3335
$ws.Range("J126").Value = 19610276
$ws.Range("K126").Value = 14002.0005
$ws.Range("L126").Value = 58830828
$ws.Range("M126").Value = -11532.0005
$ws.Range("N126").Value = -58835768
$ws.Range("H132").Value = 2565.8076
$ws.Range("I132").Value = 1712.1428
$ws.Range("J132").Value = 4323.353
$ws.Range("K132").Value = 5136.428400000001
$ws.Range("L132").Value = 12970.059
$ws.Range("M132").Value = -2606.428400000001
$ws.Range("N132").Value = -18030.059

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 12500
$ws.Range("J47").Value = 12500
$ws.Range("L47").Value = 12500
$ws.Range("N47").Value = -13480
$ws.Range("H52").Value = 12500
$ws.Range("J52").Value = 12500
$ws.Range("L52").Value = 12500
$ws.Range("N52").Value = -12966
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("H68").Value = 2841.1177
$ws.Range("I68").Value = 1457.4286
$ws.Range("J68").Value = 3809.7
$ws.Range("K68").Value = 1457.4286
$ws.Range("L68").Value = 3809.7
$ws.Range("M68").Value = -708.4286
$ws.Range("N68").Value = -5307.7
$ws.Range("H71").Value = 2841.1177
$ws.Range("I71").Value = 1457.4286
$ws.Range("J71").Value = 3809.7
$ws.Range("K71").Value = 7287.143
$ws.Range("L71").Value = 19048.5
$ws.Range("M71").Value = -3543.143
$ws.Range("N71").Value = -26536.5
$ws.Range("H93").Value = 1407
$ws.Range("I93").Value = 1465.0714
$ws.Range("J93").Value = 1000.5
$ws.Range("K93").Value = 1465.0714
$ws.Range("L93").Value = 1000.5
$ws.Range("M93").Value = -217.0714
$ws.Range("N93").Value = -3496.5
$ws.Range("H100").Value = 2118.3845
$ws.Range("I100").Value = 1546.5
$ws.Range("J100").Value = 2608.5715
$ws.Range("K100").Value = 1546.5
$ws.Range("L100").Value = 2608.5715
$ws.Range("M100").Value = -1005.5
$ws.Range("N100").Value = -3690.5715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 13996
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 13996
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 13996
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -14612
$ws.Range("H62").Value = 6949011
$ws.Range("I62").Value = 31251350
$ws.Range("J62").Value = 5485.7144
$ws.Range("K62").Value = 31251350
$ws.Range("L62").Value = 5485.7144
$ws.Range("M62").Value = -31250726
$ws.Range("N62").Value = -6733.7144
$ws.Range("H65").Value = 6949011
$ws.Range("I65").Value = 31251350
$ws.Range("J65").Value = 5485.7144
$ws.Range("K65").Value = 156256750
$ws.Range("L65").Value = 27428.572
$ws.Range("M65").Value = -156253630
$ws.Range("N65").Value = -33668.572
$ws.Range("H81").Value = 143737.22
$ws.Range("I81").Value = 125776.25
$ws.Range("J81").Value = 167685.17
$ws.Range("K81").Value = 251552.5
$ws.Range("L81").Value = 335370.34
$ws.Range("M81").Value = -250491.5
$ws.Range("N81").Value = -337492.34
$ws.Range("H84").Value = 143737.22
$ws.Range("I84").Value = 125776.25
$ws.Range("J84").Value = 167685.17
$ws.Range("K84").Value = 1257762.5
$ws.Range("L84").Value = 1676851.7
$ws.Range("M84").Value = -1252458.5
$ws.Range("N84").Value = -1687459.7
$ws.Range("H96").Value = 71429970
$ws.Range("I96").Value = 250001600
$ws.Range("J96").Value = 1316
$ws.Range("K96").Value = 250001600
$ws.Range("L96").Value = 1316
$ws.Range("M96").Value = -250000227
$ws.Range("N96").Value = -4062
